{"js": "// Rewrite the \"user story\" answer paragraph (Question 1, part 2) so that it\n// describes the expanded requirements (browser support, required sections,\n// etc.) instead of the short original description.\n//\n// The rest of the document (cover page fields, headers, etc.) is untouched\n// because those pieces live in drawing canvases / content controls that the\n// Word JS API does not expose as editable ranges.\n\nconst oldText =\n  \"As a student, I want a software that works as a web platform for discussions between students and teachers. The discussion topics will be set of programming challenge questions. I want a conversation platform compatible for uploading and downloading content so that every students can grow up in programming field by discussing several type of challenge questions.\";\n\nconst newText =\n  \"As a student, I want a software that works as a web platform which supports on most browsers like Chrome, Firefox, Safari and Microsoft edge in least, for discussions between students and teachers. It should contain login options, history section, question section, answer section and discussion section. The discussion topics will be set of programming challenge questions. The web software should be compatible for uploading and downloading content so that every students can increase their creativity in programming field by discussing several type of challenge questions.\";\n\nconst body = context.document.body;\nconst results = body.search(oldText, { matchCase: true });\nresults.load(\"items\");\nawait context.sync();\n\nif (results.items.length > 0) {\n  const target = results.items[0];\n  target.insertText(newText, Word.InsertLocation.replace);\n  await context.sync();\n}\n", "ps1": "# Rewrite the \"user story\" answer paragraph (Question 1, part 2) so that it\n# describes the expanded requirements (browser support, required sections,\n# etc.) instead of the short original description.\n#\n# The rest of the document (cover page fields, headers, etc.) is untouched\n# because those pieces live in drawing canvases / content controls that are\n# not reachable through the Word automation surface exposed here.\n\n$d = $word.ActiveDocument\n\n$oldText = \"As a student, I want a software that works as a web platform for discussions between students and teachers. The discussion topics will be set of programming challenge questions. I want a conversation platform compatible for uploading and downloading content so that every students can grow up in programming field by discussing several type of challenge questions.\"\n$newText = \"As a student, I want a software that works as a web platform which supports on most browsers like Chrome, Firefox, Safari and Microsoft edge in least, for discussions between students and teachers. It should contain login options, history section, question section, answer section and discussion section. The discussion topics will be set of programming challenge questions. The web software should be compatible for uploading and downloading content so that every students can increase their creativity in programming field by discussing several type of challenge questions.\"\n\n$rng = $d.Content\n$find = $rng.Find\n$find.Text = $oldText\n$found = $find.Execute()\n\nif ($found) {\n    $rng.Text = $newText\n}\n"}
